$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 289
$daysToAdd = 21

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $cellA.Value2 = $cellA.Value2 + $daysToAdd
    $cellB.Value2 = $cellB.Value2 + $daysToAdd
}
